$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") new values, keyed by row number
$newK = @{
    2  = 4
    3  = 1
    4  = 3
    5  = 4
    6  = 4
    7  = 1
    8  = 2
    9  = 0
    10 = 2
    11 = 1
    12 = 1
    13 = 0
    14 = 3
    15 = 3
    17 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
